# Plan Restructuracion Modulo Documentacion Masiva - apply pending edits to Hoja3
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja3")
$ws.Activate()

# --- Fill in the "ESTATUS" column (D) for rows 10-16 with "OK", matching the
# --- look of the already-completed rows above (D6:D8), i.e. the "Bueno"
# --- (green) cell style. Copy the formatting from the matching template rows
# --- so the same border treatment (incl. the thicker bottom border for the
# --- last row of the block) is preserved, then stamp in the value.

# Rows 10-15 take the formatting used by rows 6-7 (interior rows of a block).
$ws.Range("D7").Copy()
$ws.Range("D10:D15").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("D10:D15").Value = "OK"

# Row 16 is the last row of its block (thicker bottom border), matching D8.
$ws.Range("D8").Copy()
$ws.Range("D16").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("D16").Value = "OK"

# --- Update the sheet view: scroll position, zoom, and current selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$win.Zoom = 160
$ws.Range("D17").Select()

$excel.CutCopyMode = 0
